# "User Already Exists then retry the Process" — the automation re-ran its
# test-data generator for the STAGE sheet: row 2 (School/Classroom/Section
# identifiers plus the Portfolio/Assignment course columns) got fresh random
# identifiers, and the three retry rows underneath (E3:E5) got fresh
# randomized numeric-looking IDs as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

# Row 2 - School / Classroom / Section identifiers
$ws.Cells.Item(2, 1).Value = "FPK12School67577"
$ws.Cells.Item(2, 2).Value = "FPK12Classroom25206"
$ws.Cells.Item(2, 3).Value = "FPK12Section2582"

# Row 2 - Portfolio course / assignment name columns (M,N,O,P)
$ws.Cells.Item(2, 13).Value = "PortfolioCourse59437"
$ws.Cells.Item(2, 14).Value = "AssignmentName59437"
$ws.Cells.Item(2, 15).Value = "PortfolioCourse66928"
$ws.Cells.Item(2, 16).Value = "AssignmentName66928"

# Rows 3-5, column E hold numeric-looking IDs that must stay text (they were
# originally stored as shared strings, not numbers) -- force the number
# format to Text first so the new values aren't auto-coerced into numbers.
$retryIds = $ws.Range("E3:E5")
$retryIds.NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "54281"
$ws.Cells.Item(4, 5).Value = "34455"
$ws.Cells.Item(5, 5).Value = "34197"
